$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 313. This shifts the existing rows
# 313-404 down to 314-405 (and carries formatting, e.g. the date style on
# column D, along with them).
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new data record.
$ws.Cells.Item(313, 1).Value = 2
$ws.Cells.Item(313, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(313, 3).Value = "Coquimbo"
$ws.Cells.Item(313, 4).Value = 44588
$ws.Cells.Item(313, 5).Value = 4
$ws.Cells.Item(313, 6).Value = "Fruta"
$ws.Cells.Item(313, 7).Value = 100102
$ws.Cells.Item(313, 8).Value = "Cítricos"
$ws.Cells.Item(313, 9).Value = 100102005
$ws.Cells.Item(313, 10).Value = "Naranja"
$ws.Cells.Item(313, 11).Value = "Navel Late"
$ws.Cells.Item(313, 12).Value = "Primera"
$ws.Cells.Item(313, 13).Value = 16
$ws.Cells.Item(313, 14).Value = 215000
$ws.Cells.Item(313, 15).Value = 220000
$ws.Cells.Item(313, 16).Value = 217500
$ws.Cells.Item(313, 17).Value = "`$/bins (400 kilos)"
$ws.Cells.Item(313, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(313, 19).Value = 544
$ws.Cells.Item(313, 20).Value = 400
